$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 148.6
$ws.Range("J2").Value = 196
$ws.Range("L2").Value = 196
$ws.Range("N2").Value = -422

$ws.Range("H32").Value = 1899.5
$ws.Range("I32").Value = 1999.8334
$ws.Range("J32").Value = 1749
$ws.Range("K32").Value = 1999.8334
$ws.Range("L32").Value = 1749
$ws.Range("M32").Value = -1673.8334
$ws.Range("N32").Value = -2401

$ws.Range("H40").Value = 4972.227
$ws.Range("I40").Value = 2296.6667
$ws.Range("J40").Value = 5394.684
$ws.Range("K40").Value = 2296.6667
$ws.Range("L40").Value = 5394.684
$ws.Range("M40").Value = -2121.6667
$ws.Range("N40").Value = -5744.684

$ws.Range("H43").Value = 4558.9165
$ws.Range("I43").Value = 4087
$ws.Range("J43").Value = 4794.875
$ws.Range("K43").Value = 4087
$ws.Range("L43").Value = 4794.875
$ws.Range("M43").Value = -4018
$ws.Range("N43").Value = -4932.875

$ws.Range("H62").Value = 2668.3333
$ws.Range("I62").Value = 2668.3333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2668.3333
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("M62").Value = -2044.3333

$ws.Range("H65").Value = 2668.3333
$ws.Range("I65").Value = 2668.3333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13341.6665
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("M65").Value = -10221.6665

$ws.Range("H80").Value = 2024.0312
$ws.Range("I80").Value = 1174.8572
$ws.Range("J80").Value = 2684.5
$ws.Range("K80").Value = 3524.5716
$ws.Range("L80").Value = 8053.5
$ws.Range("M80").Value = -2526.5716
$ws.Range("N80").Value = -10049.5

$ws.Range("H82").Value = 2050
$ws.Range("I82").Value = 2050
$ws.Range("K82").Value = 6150
$ws.Range("M82").Value = -5744

$ws.Range("H83").Value = 2024.0312
$ws.Range("I83").Value = 1174.8572
$ws.Range("J83").Value = 2684.5
$ws.Range("K83").Value = 10573.7148
$ws.Range("L83").Value = 24160.5
$ws.Range("M83").Value = -5581.7148
$ws.Range("N83").Value = -34144.5

$ws.Range("H85").Value = 2050
$ws.Range("I85").Value = 2050
$ws.Range("K85").Value = 6150
$ws.Range("M85").Value = -4746

$ws.Range("H86").Value = 11594.866
$ws.Range("I86").Value = 11595.223
$ws.Range("J86").Value = 11594.333
$ws.Range("K86").Value = 11595.223
$ws.Range("L86").Value = 11594.333
$ws.Range("M86").Value = -10472.223
$ws.Range("N86").Value = -13840.333

$ws.Range("H87").Value = 70000
$ws.Range("J87").Value = 70000
$ws.Range("L87").Value = 70000
$ws.Range("N87").Value = -72496

$ws.Range("H89").Value = 11594.866
$ws.Range("I89").Value = 11595.223
$ws.Range("J89").Value = 11594.333
$ws.Range("K89").Value = 57976.115
$ws.Range("L89").Value = 57971.665
$ws.Range("M89").Value = -52360.115
$ws.Range("N89").Value = -69203.66500000001

$ws.Range("H90").Value = 70000
$ws.Range("J90").Value = 70000
$ws.Range("L90").Value = 210000
$ws.Range("N90").Value = -222480

$ws.Range("H98").Value = 1406.2222
$ws.Range("I98").Value = 1036.8572
$ws.Range("J98").Value = 2699
$ws.Range("K98").Value = 1036.8572
$ws.Range("L98").Value = 2699
$ws.Range("M98").Value = 461.1428000000001
$ws.Range("N98").Value = -5695

$ws.Range("H100").Value = 4239.2144
$ws.Range("I100").Value = 1869.9
$ws.Range("K100").Value = 1869.9
$ws.Range("M100").Value = -1328.9

$ws.Range("H103").Value = 700
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H112").Value = 1946.862
$ws.Range("J112").Value = 1862.2727
$ws.Range("L112").Value = 5586.8181
$ws.Range("N112").Value = -7802.8181

$ws.Range("H116").Value = 14063.263
$ws.Range("J116").Value = 14711.286
$ws.Range("L116").Value = 14711.286
$ws.Range("N116").Value = -21595.286

$ws.Range("H122").Value = 1406.2222
$ws.Range("I122").Value = 1036.8572
$ws.Range("J122").Value = 2699
$ws.Range("K122").Value = 3110.5716
$ws.Range("L122").Value = 8097
$ws.Range("M122").Value = -660.5715999999998
$ws.Range("N122").Value = -12997

$ws.Range("H137").Value = 17862722
$ws.Range("I137").Value = 20835300
$ws.Range("J137").Value = 27249.25
$ws.Range("K137").Value = 62505900
$ws.Range("L137").Value = 81747.75
$ws.Range("M137").Value = -62503350
$ws.Range("N137").Value = -86847.75

$ws.Range("H138").Value = 3197.7188
$ws.Range("I138").Value = 3102.0386
$ws.Range("J138").Value = 3612.3333
$ws.Range("K138").Value = 9306.1158
$ws.Range("L138").Value = 10836.9999
$ws.Range("M138").Value = -4166.1158
$ws.Range("N138").Value = -21116.9999

$ws.Range("H141").Value = 10159.65
$ws.Range("I141").Value = 4465.8887
$ws.Range("J141").Value = 14818.182
$ws.Range("K141").Value = 13397.6661
$ws.Range("L141").Value = 44454.546
$ws.Range("M141").Value = -8217.666100000002
$ws.Range("N141").Value = -54814.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 876
$ws.Range("I2").Value = 900.5263
$ws.Range("J2").Value = 759.5
$ws.Range("K2").Value = 900.5263
$ws.Range("L2").Value = 759.5
$ws.Range("M2").Value = -787.5263
$ws.Range("N2").Value = -985.5

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H45").Value = 1831
$ws.Range("I45").Value = 1573.6666
$ws.Range("J45").Value = 3375
$ws.Range("K45").Value = 1573.6666
$ws.Range("L45").Value = 3375
$ws.Range("M45").Value = -1196.6666
$ws.Range("N45").Value = -4129

$ws.Range("H61").Value = 5891520
$ws.Range("I61").Value = 13226.1
$ws.Range("J61").Value = 14289083
$ws.Range("K61").Value = 13226.1
$ws.Range("L61").Value = 14289083
$ws.Range("M61").Value = -13014.1
$ws.Range("N61").Value = -14289507

$ws.Range("H74").Value = 980079.5
$ws.Range("I74").Value = 1238434.1
$ws.Range("K74").Value = 1238434.1
$ws.Range("M74").Value = -1237560.1

$ws.Range("H77").Value = 980079.5
$ws.Range("I77").Value = 1238434.1
$ws.Range("K77").Value = 6192170.5
$ws.Range("M77").Value = -6187802.5

$ws.Range("H116").Value = 876
$ws.Range("I116").Value = 900.5263
$ws.Range("J116").Value = 759.5
$ws.Range("K116").Value = 900.5263
$ws.Range("L116").Value = 759.5
$ws.Range("M116").Value = 1393.4737
$ws.Range("N116").Value = -5347.5

$ws.Range("H132").Value = 6265.0684
$ws.Range("J132").Value = 7692.409
$ws.Range("L132").Value = 23077.227
$ws.Range("N132").Value = -28137.227

$ws.Range("H136").Value = 5891520
$ws.Range("I136").Value = 13226.1
$ws.Range("J136").Value = 14289083
$ws.Range("K136").Value = 39678.3
$ws.Range("L136").Value = 42867249
$ws.Range("M136").Value = -37128.3
$ws.Range("N136").Value = -42872349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 876
$ws.Range("I3").Value = 900.5263
$ws.Range("J3").Value = 759.5
$ws.Range("K3").Value = 900.5263
$ws.Range("L3").Value = 759.5
$ws.Range("M3").Value = -786.5263
$ws.Range("N3").Value = -987.5

$ws.Range("H107").Value = 1460.0244
$ws.Range("I107").Value = 1213
$ws.Range("K107").Value = 1213
$ws.Range("M107").Value = 707

$ws.Range("H134").Value = 11907897
$ws.Range("I134").Value = 4032.4443
$ws.Range("J134").Value = 33334854
$ws.Range("K134").Value = 12097.3329
$ws.Range("L134").Value = 100004562
$ws.Range("M134").Value = -9562.332900000001
$ws.Range("N134").Value = -100009632

$ws.Range("H135").Value = 62043.855
$ws.Range("J135").Value = 62043.855
$ws.Range("L135").Value = 62043.855
$ws.Range("N135").Value = -72183.85500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7549
$ws.Range("I16").Value = 8324.75
$ws.Range("K16").Value = 8324.75
$ws.Range("M16").Value = -8037.75

$ws.Range("H31").Value = 2025858.9
$ws.Range("I31").Value = 2194434
$ws.Range("J31").Value = 2957
$ws.Range("K31").Value = 2194434
$ws.Range("L31").Value = 2957
$ws.Range("M31").Value = -2194139
$ws.Range("N31").Value = -3547

$ws.Range("H34").Value = 2025858.9
$ws.Range("I34").Value = 2194434
$ws.Range("J34").Value = 2957
$ws.Range("K34").Value = 2194434
$ws.Range("L34").Value = 2957
$ws.Range("M34").Value = -2194232
$ws.Range("N34").Value = -3361

$ws.Range("H58").Value = 15040160
$ws.Range("I58").Value = 23817366
$ws.Range("J58").Value = 7360104
$ws.Range("K58").Value = 23817366
$ws.Range("L58").Value = 7360104
$ws.Range("M58").Value = -23817163
$ws.Range("N58").Value = -7360510

$ws.Range("H94").Value = 11352.6
$ws.Range("J94").Value = 1984.3334
$ws.Range("L94").Value = 1984.3334
$ws.Range("N94").Value = -2886.3334

$ws.Range("H99").Value = 16931.4
$ws.Range("I99").Value = 25817.777
$ws.Range("J99").Value = 3601.8333
$ws.Range("K99").Value = 25817.777
$ws.Range("L99").Value = 3601.8333
$ws.Range("M99").Value = -24319.777
$ws.Range("N99").Value = -6597.8333

$ws.Range("H113").Value = 7549
$ws.Range("I113").Value = 8324.75
$ws.Range("K113").Value = 8324.75
$ws.Range("M113").Value = -6154.75

$ws.Range("H122").Value = 52120
$ws.Range("J122").Value = 85866.664
$ws.Range("L122").Value = 257599.992
$ws.Range("N122").Value = -262499.992

$ws.Range("H126").Value = 16931.4
$ws.Range("I126").Value = 25817.777
$ws.Range("J126").Value = 3601.8333
$ws.Range("K126").Value = 77453.33099999999
$ws.Range("L126").Value = 10805.4999
$ws.Range("M126").Value = -74983.33099999999
$ws.Range("N126").Value = -15745.4999

$ws.Range("H132").Value = 8957.764999999999
$ws.Range("I132").Value = 10070.143
$ws.Range("J132").Value = 3766.6667
$ws.Range("K132").Value = 30210.429
$ws.Range("L132").Value = 11300.0001
$ws.Range("M132").Value = -27680.429
$ws.Range("N132").Value = -16360.0001

$ws.Range("H134").Value = 2371.1904
$ws.Range("I134").Value = 2150
$ws.Range("J134").Value = 3698.3333
$ws.Range("K134").Value = 6450
$ws.Range("L134").Value = 11094.9999
$ws.Range("M134").Value = -3915
$ws.Range("N134").Value = -16164.9999

$ws.Range("H136").Value = 15040160
$ws.Range("I136").Value = 23817366
$ws.Range("J136").Value = 7360104
$ws.Range("K136").Value = 71452098
$ws.Range("L136").Value = 22080312
$ws.Range("M136").Value = -71449548
$ws.Range("N136").Value = -22085412

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2145443.8
$ws.Range("I4").Value = 6250379
$ws.Range("K4").Value = 18751137
$ws.Range("M4").Value = -18751025

$ws.Range("H9").Value = 576690.6
$ws.Range("J9").Value = 576690.6
$ws.Range("L9").Value = 1730071.8
$ws.Range("N9").Value = -1730519.8

$ws.Range("H50").Value = 537.05554
$ws.Range("I50").Value = 632.1539
$ws.Range("J50").Value = 289.8
$ws.Range("K50").Value = 1896.4617
$ws.Range("L50").Value = 869.4000000000001
$ws.Range("M50").Value = -1415.4617
$ws.Range("N50").Value = -1831.4

$ws.Range("H53").Value = 537.05554
$ws.Range("I53").Value = 632.1539
$ws.Range("J53").Value = 289.8
$ws.Range("K53").Value = 1896.4617
$ws.Range("L53").Value = 869.4000000000001
$ws.Range("M53").Value = -1415.4617
$ws.Range("N53").Value = -1831.4

$ws.Range("H122").Value = 1076904.8
$ws.Range("J122").Value = 1815.1
$ws.Range("L122").Value = 16335.9
$ws.Range("N122").Value = -21235.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2466.3333
$ws.Range("I80").Value = 1699.5
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 1699.5
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -701.5
$ws.Range("N80").Value = -5996

$ws.Range("H83").Value = 2466.3333
$ws.Range("I83").Value = 1699.5
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 8497.5
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -3505.5
$ws.Range("N83").Value = -29984

$ws.Range("H97").Value = 2758
$ws.Range("I97").Value = 2624
$ws.Range("K97").Value = 2624
$ws.Range("M97").Value = -2128

$ws.Range("H102").Value = 12500
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 12500
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("L102").Value = 12500
$ws.Range("N102").Value = -15744

$ws.Range("H122").Value = 3759.5264
$ws.Range("I122").Value = 4819.375
$ws.Range("J122").Value = 2988.7273
$ws.Range("K122").Value = 14458.125
$ws.Range("L122").Value = 8966.1819
$ws.Range("M122").Value = -12008.125
$ws.Range("N122").Value = -13866.1819

$ws.Range("H132").Value = 52200
$ws.Range("I132").Value = 100000
$ws.Range("J132").Value = 36266.668
$ws.Range("K132").Value = 300000
$ws.Range("L132").Value = 108800.004
$ws.Range("M132").Value = -297470
$ws.Range("N132").Value = -113860.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3487.9092
$ws.Range("I22").Value = 2879.1
$ws.Range("J22").Value = 3995.25
$ws.Range("K22").Value = 2879.1
$ws.Range("L22").Value = 3995.25
$ws.Range("M22").Value = -2584.1
$ws.Range("N22").Value = -4585.25

$ws.Range("H27").Value = 3487.9092
$ws.Range("I27").Value = 2879.1
$ws.Range("J27").Value = 3995.25
$ws.Range("K27").Value = 2879.1
$ws.Range("L27").Value = 3995.25
$ws.Range("M27").Value = -2772.1
$ws.Range("N27").Value = -4209.25

$ws.Range("H40").Value = 4283.5
$ws.Range("I40").Value = 3870.5557
$ws.Range("K40").Value = 3870.5557
$ws.Range("M40").Value = -3734.5557

$ws.Range("H68").Value = 6555.1924
$ws.Range("I68").Value = 5714.647
$ws.Range("K68").Value = 5714.647
$ws.Range("M68").Value = -4965.647

$ws.Range("H71").Value = 6555.1924
$ws.Range("I71").Value = 5714.647
$ws.Range("K71").Value = 28573.235
$ws.Range("M71").Value = -24829.235

$ws.Range("H122").Value = 7400
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -28900

$ws.Range("H136").Value = 15627491
$ws.Range("I136").Value = 9616843
$ws.Range("K136").Value = 28850529
$ws.Range("M136").Value = -28847979

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 53818.914
$ws.Range("I122").Value = 5035.5293
$ws.Range("K122").Value = 15106.5879
$ws.Range("M122").Value = -12656.5879

$ws.Range("I126").Value = 3200
$ws.Range("K126").Value = 9600
$ws.Range("M126").Value = -7130

$ws.Range("H132").Value = 9262457
$ws.Range("I132").Value = 10419670
$ws.Range("K132").Value = 31259010
$ws.Range("M132").Value = -31256480
